{"js": "// Append two new list items to the end of the document body, matching the\n// list style/numbering (\"Paragraphedeliste\", ilvl 0, numId 1) used by the\n// existing bullet list, as described by the diff (semaine 6 v6.3 update).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph currently ends with \"... affichier les quartiles.\"\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert the first new bullet right after the last paragraph; it inherits\n// the same list/paragraph formatting (Paragraphedeliste, ilvl 0, numId 1).\nconst firstNew = lastParagraph.insertParagraph(\n  \"Probleme avec le http.post, il ne s'execute pas, ou pas au bon moment\u2026\",\n  Word.InsertLocation.after\n);\n\n// Insert the second new bullet right after the first new one.\nconst secondNew = firstNew.insertParagraph(\n  \"Faire un graphe avec abscisse : instance id, ordonn\u00e9e : list des status.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Append two new list items to the end of the document body, matching the\n# list style/numbering (\"Paragraphedeliste\", ilvl 0, numId 1) used by the\n# existing bullet list, as described by the diff (semaine 6 v6.3 update).\n\n$d = $word.ActiveDocument\n\n# Locate the current last paragraph (\"... affichier les quartiles.\")\n$lastIndex = $d.Paragraphs.Count\n$lastParagraph = $d.Paragraphs.Item($lastIndex)\n\n# Insert a new paragraph after it; it inherits the same paragraph/list\n# formatting (style Paragraphedeliste, ilvl 0, numId 1) as the original.\n$lastParagraph.Range.InsertParagraphAfter()\n$firstNewIndex = $d.Paragraphs.Count\n$firstNewParagraph = $d.Paragraphs.Item($firstNewIndex)\n$firstNewParagraph.Range.Text = \"Probleme avec le http.post, il ne s'execute pas, ou pas au bon moment\u2026\"\n\n# Insert the second new paragraph right after the first new one.\n$firstNewParagraph.Range.InsertParagraphAfter()\n$secondNewIndex = $d.Paragraphs.Count\n$secondNewParagraph = $d.Paragraphs.Item($secondNewIndex)\n$secondNewParagraph.Range.Text = \"Faire un graphe avec abscisse : instance id, ordonn\u00e9e : list des status.\"\n"}
